# Update loading_percent values for rows 2-25 (data rows 0-23) on columns
# B, C, E, F, G, H, I, J, N, matching the "case with 380 kV done" recompute.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 13.8401708162692
$ws.Cells.Item(2, 3).Value = 9.698000264760836
$ws.Cells.Item(2, 5).Value = 23.1620564085656
$ws.Cells.Item(2, 6).Value = 40.53173586766559
$ws.Cells.Item(2, 7).Value = 27.72711812156854
$ws.Cells.Item(2, 8).Value = 13.56646574590677
$ws.Cells.Item(2, 9).Value = 18.44498993411096
$ws.Cells.Item(2, 10).Value = 7.740156476360871
$ws.Cells.Item(2, 14).Value = 16.11008585924913
$ws.Cells.Item(3, 2).Value = 13.21169926634978
$ws.Cells.Item(3, 3).Value = 9.064147127781697
$ws.Cells.Item(3, 5).Value = 22.96171323551449
$ws.Cells.Item(3, 6).Value = 40.24183036564349
$ws.Cells.Item(3, 7).Value = 27.48693160992143
$ws.Cells.Item(3, 8).Value = 13.60311421318064
$ws.Cells.Item(3, 9).Value = 18.54396862518594
$ws.Cells.Item(3, 10).Value = 7.766173648753321
$ws.Cells.Item(3, 14).Value = 16.16446645961638
$ws.Cells.Item(4, 2).Value = 12.81212016026835
$ws.Cells.Item(4, 3).Value = 8.650560436789192
$ws.Cells.Item(4, 5).Value = 22.84270201091407
$ws.Cells.Item(4, 6).Value = 40.0757361857871
$ws.Cells.Item(4, 7).Value = 27.35290458361353
$ws.Cells.Item(4, 8).Value = 13.62930535219188
$ws.Cells.Item(4, 9).Value = 18.6104567490327
$ws.Cells.Item(4, 10).Value = 7.783365264415436
$ws.Cells.Item(4, 14).Value = 16.19978560492988
$ws.Cells.Item(5, 2).Value = 12.6460739672477
$ws.Cells.Item(5, 3).Value = 8.47586900063639
$ws.Cells.Item(5, 5).Value = 22.79525312589652
$ws.Cells.Item(5, 6).Value = 40.01110502251946
$ws.Cells.Item(5, 7).Value = 27.30172548302748
$ws.Cells.Item(5, 8).Value = 13.640901717298
$ws.Cells.Item(5, 9).Value = 18.63897992145912
$ws.Cells.Item(5, 10).Value = 7.790676892655202
$ws.Cells.Item(5, 14).Value = 16.2146645821906
$ws.Cells.Item(6, 2).Value = 12.61831531768578
$ws.Cells.Item(6, 3).Value = 8.446489618893448
$ws.Cells.Item(6, 5).Value = 22.7874388648753
$ws.Cells.Item(6, 6).Value = 40.00055909644259
$ws.Cells.Item(6, 7).Value = 27.29343633723546
$ws.Cells.Item(6, 8).Value = 13.64288291723558
$ws.Cells.Item(6, 9).Value = 18.64380222116723
$ws.Cells.Item(6, 10).Value = 7.791909455803912
$ws.Cells.Item(6, 14).Value = 16.21716461692019
$ws.Cells.Item(7, 2).Value = 12.80989349171763
$ws.Cells.Item(7, 3).Value = 8.648229415498964
$ws.Cells.Item(7, 5).Value = 22.84205779507434
$ws.Cells.Item(7, 6).Value = 40.07485211039879
$ws.Cells.Item(7, 7).Value = 27.3522003768472
$ws.Cells.Item(7, 8).Value = 13.6294580124968
$ws.Cells.Item(7, 9).Value = 18.61083564872165
$ws.Cells.Item(7, 10).Value = 7.783462633136004
$ws.Cells.Item(7, 14).Value = 16.19998429832614
$ws.Cells.Item(8, 2).Value = 13.62644396636352
$ws.Cells.Item(8, 3).Value = 9.484501512454854
$ws.Cells.Item(8, 5).Value = 23.09217436729073
$ws.Cells.Item(8, 6).Value = 40.42934196999794
$ws.Cells.Item(8, 7).Value = 27.64154723476315
$ws.Cells.Item(8, 8).Value = 13.57833398199988
$ws.Cells.Item(8, 9).Value = 18.47792667619214
$ws.Cells.Item(8, 10).Value = 7.748874548167152
$ws.Cells.Item(8, 14).Value = 16.12843648441812
$ws.Cells.Item(9, 2).Value = 15.11051454686994
$ws.Cells.Item(9, 3).Value = 10.93160956846299
$ws.Cells.Item(9, 5).Value = 23.61240379478199
$ws.Cells.Item(9, 6).Value = 41.21611051538611
$ws.Cells.Item(9, 7).Value = 28.31273278288948
$ws.Cells.Item(9, 8).Value = 13.50753372243102
$ws.Cells.Item(9, 9).Value = 18.26300937863772
$ws.Cells.Item(9, 10).Value = 7.690708659470308
$ws.Cells.Item(9, 14).Value = 16.00339357168977
$ws.Cells.Item(10, 2).Value = 16.11994314934334
$ws.Cells.Item(10, 3).Value = 11.87835212206626
$ws.Cells.Item(10, 5).Value = 24.00998395597477
$ws.Cells.Item(10, 6).Value = 41.84577234380665
$ws.Cells.Item(10, 7).Value = 28.86460350299448
$ws.Cells.Item(10, 8).Value = 13.47371164471614
$ws.Cells.Item(10, 9).Value = 18.13349014116784
$ws.Cells.Item(10, 10).Value = 7.653871302270613
$ws.Cells.Item(10, 14).Value = 15.92076636095315
$ws.Cells.Item(11, 2).Value = 16.56006267449019
$ws.Cells.Item(11, 3).Value = 12.28400482487095
$ws.Cells.Item(11, 5).Value = 24.19357746380497
$ws.Cells.Item(11, 6).Value = 42.142434481191
$ws.Cells.Item(11, 7).Value = 29.12731375131457
$ws.Cells.Item(11, 8).Value = 13.46232199622437
$ws.Cells.Item(11, 9).Value = 18.08084063683423
$ws.Cells.Item(11, 10).Value = 7.638395337577883
$ws.Cells.Item(11, 14).Value = 15.88517111870989
$ws.Cells.Item(12, 2).Value = 16.72388080451777
$ws.Cells.Item(12, 3).Value = 12.43403526981036
$ws.Cells.Item(12, 5).Value = 24.26343992702556
$ws.Cells.Item(12, 6).Value = 42.25615394212991
$ws.Cells.Item(12, 7).Value = 29.22837169382835
$ws.Cells.Item(12, 8).Value = 13.45858699592984
$ws.Cells.Item(12, 9).Value = 18.06181423942613
$ws.Cells.Item(12, 10).Value = 7.632719466544409
$ws.Cells.Item(12, 14).Value = 15.87197771182204
$ws.Cells.Item(13, 2).Value = 16.68872765358607
$ws.Cells.Item(13, 3).Value = 12.40188254001644
$ws.Cells.Item(13, 5).Value = 24.24837954369465
$ws.Cells.Item(13, 6).Value = 42.23160251887306
$ws.Cells.Item(13, 7).Value = 29.20653865717733
$ws.Cells.Item(13, 8).Value = 13.45936563739621
$ws.Cells.Item(13, 9).Value = 18.06587125765892
$ws.Cells.Item(13, 10).Value = 7.633933655214002
$ws.Cells.Item(13, 14).Value = 15.87480645133576
$ws.Cells.Item(14, 2).Value = 16.5735976485381
$ws.Cells.Item(14, 3).Value = 12.29641959848235
$ws.Cells.Item(14, 5).Value = 24.19931859110842
$ws.Cells.Item(14, 6).Value = 42.15176312195021
$ws.Cells.Item(14, 7).Value = 29.13559687453825
$ws.Cells.Item(14, 8).Value = 13.46200311306126
$ws.Cells.Item(14, 9).Value = 18.07925701624405
$ws.Cells.Item(14, 10).Value = 7.637924680267286
$ws.Cells.Item(14, 14).Value = 15.88407996673606
$ws.Cells.Item(15, 2).Value = 16.50270385933534
$ws.Cells.Item(15, 3).Value = 12.23135460444874
$ws.Cells.Item(15, 5).Value = 24.16930997861951
$ws.Cells.Item(15, 6).Value = 42.10303616101494
$ws.Cells.Item(15, 7).Value = 29.09234504623275
$ws.Cells.Item(15, 8).Value = 13.46369401212254
$ws.Cells.Item(15, 9).Value = 18.08757507162209
$ws.Cells.Item(15, 10).Value = 7.64039333883755
$ws.Cells.Item(15, 14).Value = 15.88979745260382
$ws.Cells.Item(16, 2).Value = 16.09078730274225
$ws.Cells.Item(16, 3).Value = 11.85133965811622
$ws.Cells.Item(16, 5).Value = 23.99803611850557
$ws.Cells.Item(16, 6).Value = 41.82658287507908
$ws.Cells.Item(16, 7).Value = 28.84766083679524
$ws.Cells.Item(16, 8).Value = 13.47453669299246
$ws.Cells.Item(16, 9).Value = 18.13705786202122
$ws.Cells.Item(16, 10).Value = 7.654908502534006
$ws.Cells.Item(16, 14).Value = 15.92313261885891
$ws.Cells.Item(17, 2).Value = 15.83312576652936
$ws.Cells.Item(17, 3).Value = 11.61181542476878
$ws.Cells.Item(17, 5).Value = 23.89362713138379
$ws.Cells.Item(17, 6).Value = 41.65954337688888
$ws.Cells.Item(17, 7).Value = 28.70047084166497
$ws.Cells.Item(17, 8).Value = 13.48221433407441
$ws.Cells.Item(17, 9).Value = 18.16902632403785
$ws.Cells.Item(17, 10).Value = 7.664141497477031
$ws.Cells.Item(17, 14).Value = 15.94409239304946
$ws.Cells.Item(18, 2).Value = 15.68313522858029
$ws.Cells.Item(18, 3).Value = 11.47169089422794
$ws.Cells.Item(18, 5).Value = 23.83383423986911
$ws.Cells.Item(18, 6).Value = 41.56443480475158
$ws.Cells.Item(18, 7).Value = 28.61691507243572
$ws.Cells.Item(18, 8).Value = 13.48700626686938
$ws.Cells.Item(18, 9).Value = 18.18800329335378
$ws.Cells.Item(18, 10).Value = 7.669572679981913
$ws.Cells.Item(18, 14).Value = 15.95633546728927
$ws.Cells.Item(19, 2).Value = 15.63204691324474
$ws.Cells.Item(19, 3).Value = 11.42384187834096
$ws.Cells.Item(19, 5).Value = 23.81363573688342
$ws.Cells.Item(19, 6).Value = 41.532401662536
$ws.Cells.Item(19, 7).Value = 28.58881716895038
$ws.Cells.Item(19, 8).Value = 13.48869319655439
$ws.Cells.Item(19, 9).Value = 18.19452954245561
$ws.Cells.Item(19, 10).Value = 7.671432291897681
$ws.Cells.Item(19, 14).Value = 15.96051300000857
$ws.Cells.Item(20, 2).Value = 15.86074037844753
$ws.Cells.Item(20, 3).Value = 11.63755695576726
$ws.Cells.Item(20, 5).Value = 23.90471508154969
$ws.Cells.Item(20, 6).Value = 41.67722546514504
$ws.Cells.Item(20, 7).Value = 28.716025928697
$ws.Cells.Item(20, 8).Value = 13.48135810058339
$ws.Cells.Item(20, 9).Value = 18.16556214245155
$ws.Cells.Item(20, 10).Value = 7.663146145083706
$ws.Cells.Item(20, 14).Value = 15.94184178238402
$ws.Cells.Item(21, 2).Value = 16.60749208477926
$ws.Cells.Item(21, 3).Value = 12.32749368267289
$ws.Cells.Item(21, 5).Value = 24.21372019193013
$ws.Cells.Item(21, 6).Value = 42.17517717208965
$ws.Cells.Item(21, 7).Value = 29.15639227298049
$ws.Cells.Item(21, 8).Value = 13.46121271009177
$ws.Cells.Item(21, 9).Value = 18.07530050100536
$ws.Cells.Item(21, 10).Value = 7.636747408976541
$ws.Cells.Item(21, 14).Value = 15.88134836141514
$ws.Cells.Item(22, 2).Value = 17.07891507548282
$ws.Cells.Item(22, 3).Value = 12.75753991050327
$ws.Cells.Item(22, 5).Value = 24.41762742927923
$ws.Cells.Item(22, 6).Value = 42.50861695784646
$ws.Cells.Item(22, 7).Value = 29.45332931960473
$ws.Cells.Item(22, 8).Value = 13.45141697887624
$ws.Cells.Item(22, 9).Value = 18.02162258333768
$ws.Cells.Item(22, 10).Value = 7.620570080765743
$ws.Cells.Item(22, 14).Value = 15.84347759334675
$ws.Cells.Item(23, 2).Value = 16.82885836822218
$ws.Cells.Item(23, 3).Value = 12.52991981705729
$ws.Cells.Item(23, 5).Value = 24.30863711959083
$ws.Cells.Item(23, 6).Value = 42.32995250800982
$ws.Cells.Item(23, 7).Value = 29.29404738565066
$ws.Cells.Item(23, 8).Value = 13.45633571500802
$ws.Cells.Item(23, 9).Value = 18.04978227247241
$ws.Cells.Item(23, 10).Value = 7.629105704503784
$ws.Cells.Item(23, 14).Value = 15.86353781042163
$ws.Cells.Item(24, 2).Value = 15.84826159007045
$ws.Cells.Item(24, 3).Value = 11.6259267452752
$ws.Cells.Item(24, 5).Value = 23.89970148872898
$ws.Cells.Item(24, 6).Value = 41.66922850575097
$ws.Cells.Item(24, 7).Value = 28.70899014638696
$ws.Cells.Item(24, 8).Value = 13.4817440266033
$ws.Cells.Item(24, 9).Value = 18.16712643717742
$ws.Cells.Item(24, 10).Value = 7.663595760745748
$ws.Cells.Item(24, 14).Value = 15.94285868210672
$ws.Cells.Item(25, 2).Value = 14.72265549818173
$ws.Cells.Item(25, 3).Value = 10.56065851794268
$ws.Cells.Item(25, 5).Value = 23.46876369304006
$ws.Cells.Item(25, 6).Value = 40.99387783991993
$ws.Cells.Item(25, 7).Value = 28.12048602497246
$ws.Cells.Item(25, 8).Value = 13.52350817164488
$ws.Cells.Item(25, 9).Value = 18.31620324695972
$ws.Cells.Item(25, 10).Value = 7.705409045983027
$ws.Cells.Item(25, 14).Value = 16.03559364439154
